# Reorder "Recorded By" (column G) values: when the value starts with
# "System, ", swap the first two comma-separated tokens so "System" moves
# from the first position to the second position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("System, ")) {
        $parts = $val.ToString().Split(", ")
        $first = $parts[0]
        $second = $parts[1]
        $parts[0] = $second
        $parts[1] = $first
        $newVal = [string]::Join(", ", $parts)
        $cell.Value = $newVal
    }
}
